$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Bo Giao duc article (was row 8 in the old sheet)
$ws.Range("A6").Value = "Bộ Giáo dục: Thầy cô được đàng hoàng dạy thêm"
$ws.Range("B6").Value = "https://vnexpress.net/bo-giao-duc-thay-co-duoc-dang-hoang-day-them-4785315.html"
$ws.Range("C6").Value = "Thầy cô được đàng hoàng dạy học sinh của mình ngoài nhà trường nhưng tuyệt đối không được ép buộc, theo Vụ trưởng Trung học."

# Row 7: Hanh trinh CEO Telegram (brand new article)
$ws.Range("A7").Value = "Hành trình của CEO Telegram trước khi bị bắt"
$ws.Range("B7").Value = "https://vnexpress.net/hanh-trinh-cua-ceo-telegram-truoc-khi-bi-bat-4785389.html"
$ws.Range("C7").Value = "Hơn 11 năm sau khi sáng lập Telegram, tỷ phú Pavel Durov bị bắt khi gần chạm giấc mơ một tỷ người dùng hoạt động hàng tháng trên nền tảng."

# Row 8: Co vat Hoang thanh (was row 6 in the old sheet)
$ws.Range("A8").Value = "Cổ vật Hoàng thành Thăng Long trưng bày ở TP HCM"
$ws.Range("B8").Value = "https://vnexpress.net/co-vat-hoang-thanh-thang-long-trung-bay-o-tp-hcm-4784829.html"
$ws.Range("C8").Value = "150 hiện vật, tài liệu, hình ảnh tại Khu di sản Hoàng thành Thăng Long trưng bày ở Bảo tàng TP HCM, quận 1."

# Row 9: Arsenal (was row 7 in the old sheet)
$ws.Range("A9").Value = "Arsenal đòi được món nợ từ Aston Villa"
$ws.Range("B9").Value = "https://vnexpress.net/arsenal-doi-duoc-mon-no-tu-aston-villa-4785356.html"
$ws.Range("C9").Value = "AnhTận dụng cơ hội tốt hơn kèm một chút may mắn, thầy trò Mikel Arteta đánh bại đối thủ khó chơi 2-0 ở vòng 2 Ngoại hạng Anh."

# Row 10: IS nhan trach nhiem (was row 9 in the old sheet, unchanged position content-wise just shifted by 1)
$ws.Range("A10").Value = "IS nhận trách nhiệm vụ đâm dao tại lễ hội ở Đức"
$ws.Range("B10").Value = "https://vnexpress.net/is-nhan-trach-nhiem-vu-dam-dao-tai-le-hoi-o-duc-4785362.html"
$ws.Range("C10").Value = "IS tuyên bố kẻ thực hiện vụ đâm dao khiến ba người chết tại lễ hội ở thành phố Solingen, phía tây Đức, là thành viên của nhóm này."

# Row 11: iPhone 16 (was row 10 in the old sheet)
$ws.Range("A11").Value = "iPhone 16 có thể bán tại Việt Nam cuối tháng 9"
$ws.Range("B11").Value = "https://vnexpress.net/iphone-16-co-the-ban-tai-viet-nam-cuoi-thang-9-4784978.html"
$ws.Range("C11").Value = "Việt Nam vẫn nằm trong những thị trường ưu tiên thứ hai của Apple và nhiều khả năng sẽ bán iPhone 16 ngay trong tháng 9."

# Remove the old rows 12-21 entirely (they are no longer part of the data)
$ws.Range("A12:C21").ClearContents()
